# Split "natural gas nonpeaker" into two distinct electricity sources:
# "natural gas steam turbine" and "natural gas combined cycle" on the
# CPPbES (CO2 Capture Potential by Electricity Source) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbES")

# Insert a new row right after the current "natural gas nonpeaker" row (row 3)
# so we end up with two rows in its place instead of one.
$ws.Rows.Item(4).Insert()

# Row 3: natural gas steam turbine (keeps the capturable=1 flag)
$ws.Range("A3").Value = "natural gas steam turbine"
$ws.Range("B3").Value = 1

# Row 4: natural gas combined cycle (also capturable=1)
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("B4").Value = 1

Write-Host "Done splitting natural gas nonpeaker into steam turbine / combined cycle rows."
